# Applies the "Automatic update of files" edit:
#  1. Updates the "Förändrad" (column C) date for every data row (2-261)
#     from 45184 (2023-09-15) to 45186 (2023-09-17).
#  2. Adds a friendly display-name second argument (the value of column A,
#     e.g. "A 67854-2019") to every HYPERLINK() formula in columns
#     S, T, V, W, X, Y for the rows that have them (rows 2-9).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bulk-update column C (Förändrad) for all data rows (rows 2-261) --
$ws.Range("C2:C261").Value = 45186

# --- 2. Rebuild the HYPERLINK formulas that get a friendly display name --
$linkCols = @(
    @{ Col = "S"; Folder = "artfynd";         Ext = "xlsx" },
    @{ Col = "T"; Folder = "kartor";           Ext = "png"  },
    @{ Col = "V"; Folder = "klagomål";        Ext = "docx" },
    @{ Col = "W"; Folder = "klagomålsmail";   Ext = "docx" },
    @{ Col = "X"; Folder = "tillsyn";          Ext = "docx" },
    @{ Col = "Y"; Folder = "tillsynsmail";    Ext = "docx" }
)

for ($row = 2; $row -le 9; $row++) {
    $caseId = $ws.Range("A$row").Value2

    foreach ($link in $linkCols) {
        $col = $link.Col
        $url = "https://klasma.github.io/Logging_ARBOGA/$($link.Folder)/$caseId.$($link.Ext)"
        $formula = '=HYPERLINK("' + $url + '", "' + $caseId + '")'
        $ws.Range("$col$row").Formula = $formula
    }
}
